$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 99,3
$arr[0,0] = 220.00026959066784
$arr[0,1] = 250.00072582102888
$arr[0,2] = 10.864927139321173
$arr[1,0] = 219.56226959066782
$arr[1,1] = 250.04219582102888
$arr[1,2] = 20.000230880600615
$arr[2,0] = 219.55563831273267
$arr[2,1] = 250.0423588955498
$arr[2,2] = 20.00000000000059
$arr[3,0] = 219.55563831273267
$arr[3,1] = 250.0423588955498
$arr[3,2] = 20.00000000000059
$arr[4,0] = 219.45563831273267
$arr[4,1] = 248.0423588955498
$arr[4,2] = 20.00000000000059
$arr[5,0] = 218.45563831273267
$arr[5,1] = 248.0423588955498
$arr[5,2] = 20.00000000000059
$arr[6,0] = 218.45563831273267
$arr[6,1] = 248.0423588955498
$arr[6,2] = 20.00000000000059
$arr[7,0] = 218.45563831273267
$arr[7,1] = 248.0423588955498
$arr[7,2] = 20.00000000000059
$arr[8,0] = 218.45563831273267
$arr[8,1] = 248.0423588955498
$arr[8,2] = 20.00000000000059
$arr[9,0] = 217.45563831273267
$arr[9,1] = 248.0423588955498
$arr[9,2] = 20.00000000000059
$arr[10,0] = 217.45563831273267
$arr[10,1] = 248.0423588955498
$arr[10,2] = 20.00000000000059
$arr[11,0] = 217.45563831273267
$arr[11,1] = 248.0423588955498
$arr[11,2] = 20.00000000000059
$arr[12,0] = 217.45563831273267
$arr[12,1] = 248.0423588955498
$arr[12,2] = 20.00000000000059
$arr[13,0] = 217.45563831273267
$arr[13,1] = 248.0423588955498
$arr[13,2] = 20.00000000000059
$arr[14,0] = 217.45563831273267
$arr[14,1] = 248.0423588955498
$arr[14,2] = 20.00000000000059
$arr[15,0] = 217.45563831273267
$arr[15,1] = 248.0423588955498
$arr[15,2] = 20.00000000000059
$arr[16,0] = 217.45563831273267
$arr[16,1] = 248.0423588955498
$arr[16,2] = 20.00000000000059
$arr[17,0] = 217.45563831273267
$arr[17,1] = 248.0423588955498
$arr[17,2] = 20.00000000000059
$arr[18,0] = 217.45563831273267
$arr[18,1] = 248.0423588955498
$arr[18,2] = 20.00000000000059
$arr[19,0] = 217.45563831273267
$arr[19,1] = 248.0423588955498
$arr[19,2] = 20.00000000000059
$arr[20,0] = 217.45563831273267
$arr[20,1] = 248.0423588955498
$arr[20,2] = 20.00000000000059
$arr[21,0] = 217.45563831273267
$arr[21,1] = 248.0423588955498
$arr[21,2] = 20.00000000000059
$arr[22,0] = 217.45563831273267
$arr[22,1] = 248.0423588955498
$arr[22,2] = 20.00000000000059
$arr[23,0] = 217.45563831273267
$arr[23,1] = 248.0423588955498
$arr[23,2] = 20.00000000000059
$arr[24,0] = 212.45563831273267
$arr[24,1] = 248.0423588955498
$arr[24,2] = 20.00000000000059
$arr[25,0] = 212.45563831273267
$arr[25,1] = 248.0423588955498
$arr[25,2] = 20.00000000000059
$arr[26,0] = 212.45563831273267
$arr[26,1] = 248.0423588955498
$arr[26,2] = 20.00000000000059
$arr[27,0] = 212.45563831273267
$arr[27,1] = 248.0423588955498
$arr[27,2] = 20.00000000000059
$arr[28,0] = 212.45563831273267
$arr[28,1] = 248.0423588955498
$arr[28,2] = 20.00000000000059
$arr[29,0] = 212.45563831273267
$arr[29,1] = 248.0423588955498
$arr[29,2] = 20.00000000000059
$arr[30,0] = 212.45563831273267
$arr[30,1] = 248.0423588955498
$arr[30,2] = 20.00000000000059
$arr[31,0] = 212.45563831273267
$arr[31,1] = 248.0423588955498
$arr[31,2] = 20.00000000000059
$arr[32,0] = 212.45563831273267
$arr[32,1] = 248.0423588955498
$arr[32,2] = 20.00000000000059
$arr[33,0] = 212.45563831273267
$arr[33,1] = 248.0423588955498
$arr[33,2] = 20.00000000000059
$arr[34,0] = 212.45563831273267
$arr[34,1] = 248.0423588955498
$arr[34,2] = 20.00000000000059
$arr[35,0] = 212.45563831273267
$arr[35,1] = 248.0423588955498
$arr[35,2] = 20.00000000000059
$arr[36,0] = 212.45563831273267
$arr[36,1] = 248.0423588955498
$arr[36,2] = 20.00000000000059
$arr[37,0] = 212.45563831273267
$arr[37,1] = 248.0423588955498
$arr[37,2] = 20.00000000000059
$arr[38,0] = 212.45563831273267
$arr[38,1] = 248.0423588955498
$arr[38,2] = 20.00000000000059
$arr[39,0] = 212.45563831273267
$arr[39,1] = 248.0423588955498
$arr[39,2] = 20.00000000000059
$arr[40,0] = 212.45563831273267
$arr[40,1] = 248.0423588955498
$arr[40,2] = 20.00000000000059
$arr[41,0] = 212.45563831273267
$arr[41,1] = 248.0423588955498
$arr[41,2] = 20.00000000000059
$arr[42,0] = 212.45563831273267
$arr[42,1] = 248.0423588955498
$arr[42,2] = 20.00000000000059
$arr[43,0] = 212.45563831273267
$arr[43,1] = 248.0423588955498
$arr[43,2] = 20.00000000000059
$arr[44,0] = 212.45563831273267
$arr[44,1] = 248.0423588955498
$arr[44,2] = 20.00000000000059
$arr[45,0] = 212.45563831273267
$arr[45,1] = 248.0423588955498
$arr[45,2] = 20.00000000000059
$arr[46,0] = 212.45563831273267
$arr[46,1] = 248.0423588955498
$arr[46,2] = 20.00000000000059
$arr[47,0] = 212.45563831273267
$arr[47,1] = 248.0423588955498
$arr[47,2] = 20.00000000000059
$arr[48,0] = 212.45563831273267
$arr[48,1] = 248.0423588955498
$arr[48,2] = 20.00000000000059
$arr[49,0] = 212.45563831273267
$arr[49,1] = 248.0423588955498
$arr[49,2] = 20.00000000000059
$arr[50,0] = 212.45563831273267
$arr[50,1] = 248.0423588955498
$arr[50,2] = 20.00000000000059
$arr[51,0] = 212.45563831273267
$arr[51,1] = 248.0423588955498
$arr[51,2] = 20.00000000000059
$arr[52,0] = 212.45563831273267
$arr[52,1] = 248.0423588955498
$arr[52,2] = 20.00000000000059
$arr[53,0] = 207.45563831273267
$arr[53,1] = 248.0423588955498
$arr[53,2] = 20.00000000000059
$arr[54,0] = 207.45563831273267
$arr[54,1] = 248.0423588955498
$arr[54,2] = 20.00000000000059
$arr[55,0] = 207.45563831273267
$arr[55,1] = 248.0423588955498
$arr[55,2] = 20.00000000000059
$arr[56,0] = 207.45563831273267
$arr[56,1] = 248.0423588955498
$arr[56,2] = 20.00000000000059
$arr[57,0] = 207.45563831273267
$arr[57,1] = 248.0423588955498
$arr[57,2] = 20.00000000000059
$arr[58,0] = 207.45563831273267
$arr[58,1] = 248.0423588955498
$arr[58,2] = 20.00000000000059
$arr[59,0] = 207.45563831273267
$arr[59,1] = 248.0423588955498
$arr[59,2] = 20.00000000000059
$arr[60,0] = 207.45563831273267
$arr[60,1] = 248.0423588955498
$arr[60,2] = 20.00000000000059
$arr[61,0] = 207.45563831273267
$arr[61,1] = 248.0423588955498
$arr[61,2] = 20.00000000000059
$arr[62,0] = 207.45563831273267
$arr[62,1] = 248.0423588955498
$arr[62,2] = 20.00000000000059
$arr[63,0] = 207.45563831273267
$arr[63,1] = 248.0423588955498
$arr[63,2] = 20.00000000000059
$arr[64,0] = 207.45563831273267
$arr[64,1] = 248.0423588955498
$arr[64,2] = 20.00000000000059
$arr[65,0] = 207.45563831273267
$arr[65,1] = 248.0423588955498
$arr[65,2] = 20.00000000000059
$arr[66,0] = 207.45563831273267
$arr[66,1] = 248.0423588955498
$arr[66,2] = 20.00000000000059
$arr[67,0] = 207.45563831273267
$arr[67,1] = 248.0423588955498
$arr[67,2] = 20.00000000000059
$arr[68,0] = 207.45563831273267
$arr[68,1] = 248.0423588955498
$arr[68,2] = 20.00000000000059
$arr[69,0] = 207.45563831273267
$arr[69,1] = 248.0423588955498
$arr[69,2] = 20.00000000000059
$arr[70,0] = 207.45563831273267
$arr[70,1] = 248.0423588955498
$arr[70,2] = 20.00000000000059
$arr[71,0] = 207.45563831273267
$arr[71,1] = 248.0423588955498
$arr[71,2] = 20.00000000000059
$arr[72,0] = 207.45563831273267
$arr[72,1] = 248.0423588955498
$arr[72,2] = 20.00000000000059
$arr[73,0] = 207.45563831273267
$arr[73,1] = 248.0423588955498
$arr[73,2] = 20.00000000000059
$arr[74,0] = 207.45563831273267
$arr[74,1] = 248.0423588955498
$arr[74,2] = 20.00000000000059
$arr[75,0] = 207.45563831273267
$arr[75,1] = 248.0423588955498
$arr[75,2] = 20.00000000000059
$arr[76,0] = 207.45563831273267
$arr[76,1] = 248.0423588955498
$arr[76,2] = 20.00000000000059
$arr[77,0] = 207.45563831273267
$arr[77,1] = 248.0423588955498
$arr[77,2] = 20.00000000000059
$arr[78,0] = 207.45563831273267
$arr[78,1] = 248.0423588955498
$arr[78,2] = 20.00000000000059
$arr[79,0] = 207.45563831273267
$arr[79,1] = 248.0423588955498
$arr[79,2] = 20.00000000000059
$arr[80,0] = 207.45563831273267
$arr[80,1] = 248.0423588955498
$arr[80,2] = 20.00000000000059
$arr[81,0] = 207.45563831273267
$arr[81,1] = 248.0423588955498
$arr[81,2] = 20.00000000000059
$arr[82,0] = 207.45563831273267
$arr[82,1] = 248.0423588955498
$arr[82,2] = 20.00000000000059
$arr[83,0] = 207.45563831273267
$arr[83,1] = 248.0423588955498
$arr[83,2] = 20.00000000000059
$arr[84,0] = 207.45563831273267
$arr[84,1] = 248.0423588955498
$arr[84,2] = 20.00000000000059
$arr[85,0] = 207.45563831273267
$arr[85,1] = 248.0423588955498
$arr[85,2] = 20.00000000000059
$arr[86,0] = 207.45563831273267
$arr[86,1] = 248.0423588955498
$arr[86,2] = 20.00000000000059
$arr[87,0] = 207.45563831273267
$arr[87,1] = 248.0423588955498
$arr[87,2] = 20.00000000000059
$arr[88,0] = 207.45563831273267
$arr[88,1] = 248.0423588955498
$arr[88,2] = 20.00000000000059
$arr[89,0] = 207.45563831273267
$arr[89,1] = 248.0423588955498
$arr[89,2] = 20.00000000000059
$arr[90,0] = 207.45563831273267
$arr[90,1] = 248.0423588955498
$arr[90,2] = 20.00000000000059
$arr[91,0] = 207.45563831273267
$arr[91,1] = 248.0423588955498
$arr[91,2] = 20.00000000000059
$arr[92,0] = 207.45563831273267
$arr[92,1] = 248.0423588955498
$arr[92,2] = 20.00000000000059
$arr[93,0] = 207.45563831273267
$arr[93,1] = 248.0423588955498
$arr[93,2] = 20.00000000000059
$arr[94,0] = 207.45563831273267
$arr[94,1] = 248.0423588955498
$arr[94,2] = 20.00000000000059
$arr[95,0] = 207.45563831273267
$arr[95,1] = 248.0423588955498
$arr[95,2] = 20.00000000000059
$arr[96,0] = 207.45563831273267
$arr[96,1] = 248.0423588955498
$arr[96,2] = 20.00000000000059
$arr[97,0] = 207.45563831273267
$arr[97,1] = 248.0423588955498
$arr[97,2] = 20.00000000000059
$arr[98,0] = 207.45563831273267
$arr[98,1] = 248.0423588955498
$arr[98,2] = 20.00000000000059
$ws.Range("A3:C101").Value = $arr
